# Update "想去人数" (interest counts) on two sheets to reflect newly
# generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1182
$wsExhibit.Range("F4").Value = 2648

# Sheet "全部类型" (All types) - combined listing, same events appear later
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1182
$wsAll.Range("F6").Value = 2648
